$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at I, shifting the existing "Logistic percentage"
# column (and its data) over to J.
$ws.Columns("I:I").Insert()

# New "Distribution channel code" column (now column I).
$ws.Range("I1").Value = "Distribution channel code"
$ws.Range("I2").Value = "TR"
$ws.Range("I3").Value = "GO"

# Match the column width Excel computed for the new column's content.
$ws.Columns("I:I").ColumnWidth = 21.6
